$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.298.32"
Set-TextValue "E2" "  -0.50%  "

# Row 3
Set-TextValue "D3" "1.873.40"
Set-TextValue "E3" "  -0.19%  "

# Row 4
Set-TextValue "E4" "  +0.12%  "

# Row 5
Set-TextValue "D5" "0.7081"
Set-TextValue "E5" "  -0.84%  "

# Row 6
Set-TextValue "D6" "241.79"
Set-TextValue "E6" "  -0.13%  "

# Row 7
Set-TextValue "D7" "1.001"
Set-TextValue "E7" "  +0.08%  "

# Row 8
Set-TextValue "D8" "0.07769"
Set-TextValue "E8" "  +0.35%  "

# Row 9
Set-TextValue "D9" "0.3100"
Set-TextValue "E9" "  -0.69%  "

# Row 10
Set-TextValue "D10" "25.00"
Set-TextValue "E10" "  -0.31%  "

# Row 11
Set-TextValue "D11" "0.08387"
Set-TextValue "E11" "  +0.06%  "

# Row 12
Set-TextValue "D12" "1.880.02"
Set-TextValue "E12" "  -0.18%  "

# Row 13
Set-TextValue "D13" "5.235"
Set-TextValue "E13" "  -0.34%  "

# Row 14
Set-TextValue "D14" "0.7158"
Set-TextValue "E14" "  -0.41%  "

# Row 15
Set-TextValue "D15" "91.12"
Set-TextValue "E15" "  -0.56%  "

# Row 16
Set-TextValue "D16" "29.313.95"
Set-TextValue "E16" "  -0.42%  "

# Row 17
Set-TextValue "D17" "6.092"
Set-TextValue "E17" "  +1.71%  "

# Row 18
Set-TextValue "D18" "0.000008272"
Set-TextValue "E18" "  +0.72%  "

# Row 19
Set-TextValue "E19" "  -1.88%  "

# Row 20
Set-TextValue "D20" "13.20"
Set-TextValue "E20" "  -0.23%  "

# Row 21
Set-TextValue "D21" "2.123.78"
Set-TextValue "E21" "  -0.95%  "

# Row 22
Set-TextValue "D22" "0.9999"

# Row 23
Set-TextValue "D23" "7.752"
Set-TextValue "E23" "  -2.50%  "

# Row 24
Set-TextValue "D24" "1.001"
Set-TextValue "E24" "  +0.12%  "

# Row 25
Set-TextValue "D25" "0.1586"
Set-TextValue "E25" "  -2.74%  "

# Row 26
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "162.28"
Set-TextValue "E26" "  -1.07%  "

# Row 27
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "9.029"
Set-TextValue "E27" "  -0.23%  "

# Row 28
Set-TextValue "D28" "18.52"
Set-TextValue "E28" "  -0.46%  "

# Row 29
Set-TextValue "E29" "  -0.36%  "

# Row 30
Set-TextValue "D30" "4.402"
Set-TextValue "E30" "  -0.62%  "

# Row 31
Set-TextValue "D31" "4.307"
Set-TextValue "E31" "  -0.70%  "

# Row 32
Set-TextValue "D32" "1.272"
Set-TextValue "E32" "  -1.95%  "

# Row 33
Set-TextValue "D33" "0.05356"
Set-TextValue "E33" "  +1.99%  "

# Row 34
Set-TextValue "D34" "1.937"
Set-TextValue "E34" "  +0.19%  "

# Row 35
Set-TextValue "D35" "0.7509"
Set-TextValue "E35" "  -2.64%  "

# Row 36
Set-TextValue "D36" "1.175"
Set-TextValue "E36" "  -0.08%  "

# Row 37
Set-TextValue "D37" "2.684"
Set-TextValue "E37" "  +0.24%  "

# Row 38
Set-TextValue "D38" "0.01874"
Set-TextValue "E38" "  +0.33%  "

# Row 39
Set-TextValue "D39" "1.238.75"
Set-TextValue "E39" "  +5.63%  "

# Row 40
Set-TextValue "D40" "2.728"
Set-TextValue "E40" "  +0.23%  "

# Row 41
Set-TextValue "D41" "6.483"
Set-TextValue "E41" "  +0.73%  "

# Row 42
Set-TextValue "D42" "0.8931"
Set-TextValue "E42" "  +0.12%  "

# Row 43
Set-TextValue "D43" "72.27"
Set-TextValue "E43" "  -1.84%  "

# Row 44
Set-TextValue "D44" "108.66"
Set-TextValue "E44" "  +4.33%  "

# Row 45
Set-TextValue "E45" "  +0.11%  "

# Row 46
Set-TextValue "D46" "2.022.76"
Set-TextValue "E46" "  -0.23%  "

# Row 47
Set-TextValue "E47" "  +8.21%  "

# Row 48
Set-TextValue "D48" "0.5201"
Set-TextValue "E48" "  +0.01%  "

# Row 49
Set-TextValue "D49" "1.792"
Set-TextValue "E49" "  -1.00%  "

# Row 50
Set-TextValue "D50" "9.431"
Set-TextValue "E50" "  -0.02%  "

# Row 51
Set-TextValue "D51" "0.4335"
Set-TextValue "E51" "  +0.27%  "
